$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the species-specific data of rows 10, 11 and 12:
#   new row 10 <- old row 11
#   new row 11 <- old row 12
#   new row 12 <- old row 10
# Only columns A, B, D, E, F, G, H, Q, R and AC (the "Publik kommentar") change.
# Capture the original values first, since we will overwrite them in place.

$A10 = $ws.Range("A10").Value2
$B10 = $ws.Range("B10").Value2
$D10 = $ws.Range("D10").Value2
$E10 = $ws.Range("E10").Value2
$F10 = $ws.Range("F10").Value2
$G10 = $ws.Range("G10").Value2
$H10 = $ws.Range("H10").Value2
$Q10 = $ws.Range("Q10").Value2
$R10 = $ws.Range("R10").Value2
$AC10 = $ws.Range("AC10").Value2

$A11 = $ws.Range("A11").Value2
$B11 = $ws.Range("B11").Value2
$D11 = $ws.Range("D11").Value2
$E11 = $ws.Range("E11").Value2
$F11 = $ws.Range("F11").Value2
$G11 = $ws.Range("G11").Value2
$H11 = $ws.Range("H11").Value2
$Q11 = $ws.Range("Q11").Value2
$R11 = $ws.Range("R11").Value2

$A12 = $ws.Range("A12").Value2
$B12 = $ws.Range("B12").Value2
$D12 = $ws.Range("D12").Value2
$E12 = $ws.Range("E12").Value2
$F12 = $ws.Range("F12").Value2
$G12 = $ws.Range("G12").Value2
$H12 = $ws.Range("H12").Value2
$Q12 = $ws.Range("Q12").Value2
$R12 = $ws.Range("R12").Value2
$AC12 = $ws.Range("AC12").Value2

# Row 10 becomes old row 11 (no public comment)
$ws.Range("A10").Value = $A11
$ws.Range("B10").Value = $B11
$ws.Range("D10").Value = $D11
$ws.Range("E10").Value = $E11
$ws.Range("F10").Value = $F11
$ws.Range("G10").Value = $G11
$ws.Range("H10").Value = $H11
$ws.Range("Q10").Value = $Q11
$ws.Range("R10").Value = $R11
$ws.Range("AC10").ClearContents()

# Row 11 becomes old row 12 (gains the public comment that used to be on row 12)
$ws.Range("A11").Value = $A12
$ws.Range("B11").Value = $B12
$ws.Range("D11").Value = $D12
$ws.Range("E11").Value = $E12
$ws.Range("F11").Value = $F12
$ws.Range("G11").Value = $G12
$ws.Range("H11").Value = $H12
$ws.Range("Q11").Value = $Q12
$ws.Range("R11").Value = $R12
$ws.Range("AC11").Value = $AC12

# Row 12 becomes old row 10 (gains the public comment that used to be on row 10)
$ws.Range("A12").Value = $A10
$ws.Range("B12").Value = $B10
$ws.Range("D12").Value = $D10
$ws.Range("E12").Value = $E10
$ws.Range("F12").Value = $F10
$ws.Range("G12").Value = $G10
$ws.Range("H12").Value = $H10
$ws.Range("Q12").Value = $Q10
$ws.Range("R12").Value = $R10
$ws.Range("AC12").Value = $AC10
